$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp string
$ws.Range("A1").Value = "Datos actualizados a 26 de Marzo de 2020 a las 23:12"

# Update Cataluña row (row 5) figures
$ws.Range("B5").Value = 12940
$ws.Range("C5").Value = 2384
$ws.Range("D5").Value = 9676
$ws.Range("E5").Value = 880
